# Repull data: update column F (dSF) values per refreshed source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    4  = -3
    5  = 4
    6  = -3
    8  = -2
    9  = -2
    10 = 9
    11 = 1
    12 = -1
    13 = 3
    14 = -5
    15 = 0
    16 = 11
    17 = 3
    18 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
